$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (Total) sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalAnchor = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalAnchor)
$newSheet.Name = "2022-Q1"

# NOTE: re-fetch the "总计" worksheet *after* the insert - a reference
# obtained beforehand tracks the sheet's position, which is now occupied
# by the freshly inserted sheet, not the original "总计" worksheet.
$total = $wb.Worksheets.Item("总计")

# Clone the cell formatting (bold/bordered header row + bold row-index
# column) from the neighbouring "2021-Q4" sheet so the new sheet matches
# the look of the other quarterly sheets.
$q4.Range("A1:H8").Copy()
$newSheet.Range("A1:H8").PasteSpecial(-4122)

# ---------------------------- header row ------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------- data rows -------------------------------
# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "160910", "大成创新成长混合(LOF)",       "16.65", "85.97", "8.89", "1.4802", 1),
    @(1, "161605", "融通蓝筹成长混合",             "4.82",  "71.70", "5.43", "0.2617", 2),
    @(2, "000717", "融通转型三动力灵活配置混合A",   "3.83",  "94.89", "6.56", "0.2512", 2),
    @(3, "003165", "鹏华弘嘉灵活配置混合A",         "1.53",  "93.95", "3.29", "0.0503", 6),
    @(4, "009828", "融通转型三动力灵活配置混合C",   "0.59",  "94.89", "6.56", "0.0387", 2),
    @(5, "007903", "长城量化小盘股票",             "1.36",  "90.57", "1.42", "0.0193", 1),
    @(6, "003166", "鹏华弘嘉灵活配置混合C",         "0.56",  "93.95", "3.29", "0.0184", 6)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]

    # Columns B, D, E, F, G hold numeric-looking text (fund codes / ratios
    # that must keep trailing zeros & leading zeros), so force them to
    # text with a leading apostrophe and then clear the resulting
    # "quote prefix" style so the cell ends up with no style override -
    # matching the plain inline-string cells produced by the source tool.
    $newSheet.Range("B$r").Value = "'" + $row[1]
    $newSheet.Range("B$r").ClearFormats()

    $newSheet.Range("C$r").Value = $row[2]

    $newSheet.Range("D$r").Value = "'" + $row[3]
    $newSheet.Range("D$r").ClearFormats()

    $newSheet.Range("E$r").Value = "'" + $row[4]
    $newSheet.Range("E$r").ClearFormats()

    $newSheet.Range("F$r").Value = "'" + $row[5]
    $newSheet.Range("F$r").ClearFormats()

    $newSheet.Range("G$r").Value = "'" + $row[6]
    $newSheet.Range("G$r").ClearFormats()

    $newSheet.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: add a new 2022-Q1 summary row at the top
#    of the data (row 2) and push the existing rows down.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Give the new row's index cell (A2) the same bold style as the other
# index cells in column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 2.12

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
